$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns B (Date) and F (LR) hold text-like values that could otherwise be
# auto-converted to dates/numbers by Excel. Pre-format them as Text so the
# values entered below are preserved verbatim as strings.
$ws.Range("B2:B4").NumberFormat = "@"
$ws.Range("F2:F4").NumberFormat = "@"

# --- Row 2: corrected values for the existing LR entry ---
$ws.Range("B2").Value = "2021-01-03"
$ws.Range("C2").Value = "SILIG"
$ws.Range("D2").Value = "SSWW"
$ws.Range("E2").Value = "DARJ"
$ws.Range("F2").Value = "13"
$ws.Range("H2").Value = 111222
$ws.Range("I2").Value = 123
$ws.Range("J2").Value = 12
$ws.Range("K2").Value = 12212
$ws.Range("L2").Value = 122

# --- Row 3: new LR party route ---
$ws.Range("A3").Value = 2
$ws.Range("B3").Value = "2021-01-20"
$ws.Range("C3").Value = "A"
$ws.Range("D3").Value = "A"
$ws.Range("E3").Value = "A"
$ws.Range("F3").Value = "A"
$ws.Range("G3").Value = 121
$ws.Range("H3").Value = 123

# --- Row 4: new LR party route ---
$ws.Range("A4").Value = 3
$ws.Range("B4").Value = "2011-08-19"
$ws.Range("C4").Value = "Q"
$ws.Range("D4").Value = "Q"
$ws.Range("E4").Value = "Q"
$ws.Range("F4").Value = "Q"
$ws.Range("G4").Value = 112
$ws.Range("H4").Value = 12
